$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh (GitHub Actions cron). Cells that hold
# thousands-dotted / trailing-zero price strings look numeric to Excel's
# smart-typing, so force text via NumberFormat='@' and restore the default
# 'Normal' style afterwards (keeps cell styling identical to the original,
# unstyled inline-string cells).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "61.494.53"
$ws.Range("E2").Value = "  -4.20%  "
$ws.Range("D3").Value = "2.968.07"
$ws.Range("E3").Value = "  -5.46%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "538.80"
$ws.Range("E5").Value = "  -5.58%  "
Set-TextValue "D6" "150.36"
$ws.Range("E6").Value = "  -7.56%  "
$ws.Range("E7").Value = "  +0.14%  "
Set-TextValue "D8" "0.567"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "2.978.74"
$ws.Range("E9").Value = "  -5.40%  "
$ws.Range("E10").Value = "  -3.26%  "
Set-TextValue "D11" "6.14"
$ws.Range("E11").Value = "  -6.61%  "
Set-TextValue "D12" "0.367"
$ws.Range("E12").Value = "  -4.45%  "
$ws.Range("D13").Value = "3.492.55"
$ws.Range("E13").Value = "  -5.23%  "
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "61.575.06"
$ws.Range("E15").Value = "  -4.13%  "
Set-TextValue "D16" "23.60"
$ws.Range("E16").Value = "  -5.85%  "
$ws.Range("D17").Value = "2.974.53"
$ws.Range("E17").Value = "  -5.22%  "
$ws.Range("E18").Value = "  -4.61%  "
Set-TextValue "D19" "5.16"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  -3.94%  "
Set-TextValue "D21" "380.26"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("E22").Value = "  -6.00%  "
$ws.Range("E23").Value = "  -0.05%  "
Set-TextValue "D24" "5.65"
$ws.Range("E24").Value = "  -3.56%  "
Set-TextValue "D25" "65.41"
$ws.Range("E25").Value = "  -3.85%  "
Set-TextValue "D26" "0.469"
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").Value = "3.096.43"
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  -6.59%  "
Set-TextValue "D31" "8.21"
$ws.Range("E31").Value = "  -6.29%  "
$ws.Range("E32").Value = "  +0.02%  "
Set-TextValue "D33" "1.71"
$ws.Range("E33").Value = "  -4.79%  "
Set-TextValue "D34" "20.44"
$ws.Range("E34").Value = "  -3.29%  "
Set-TextValue "D35" "160.80"
$ws.Range("E35").Value = "  +0.67%  "
Set-TextValue "D36" "4.64"
$ws.Range("E36").Value = "  -3.47%  "
Set-TextValue "D37" "5.91"
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("E40").Value = "  -7.11%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D41" "37.51"
$ws.Range("E41").Value = "  -2.03%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D42" "3.90"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("D43").Value = "2.406.19"
$ws.Range("E43").Value = "  -9.10%  "
Set-TextValue "D44" "22.13"
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("E45").Value = "  -3.11%  "
Set-TextValue "D46" "0.0590"
$ws.Range("E46").Value = "  -3.28%  "
Set-TextValue "D47" "5.12"
$ws.Range("E47").Value = "  -5.71%  "
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  -2.92%  "
Set-TextValue "D50" "0.0950"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D51" "19.65"
$ws.Range("E51").Value = "  -6.42%  "
